$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the latest scrape.
# Column D (Price) values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "1.021", "0.00001037") are stored as text,
# matching the inline-string cell type used throughout the sheet; the Text
# number-format is cleared again immediately after so no extra formatting
# is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.916.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.021"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.020"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4635"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3892"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.80"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07875"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.901.47"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.897"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.036"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.025"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06754"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.51"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001037"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.017"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.922.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.443"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.83"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.355"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.121.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.84"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.045"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.370"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.89"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09456"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9513"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.678"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.287"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.335"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06060"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02219"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.208"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.031"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5887"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1867"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.272"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5603"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.04"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.392"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.893"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06895"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.058"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.16%  "
